$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure Price column cells keep exact text representation
foreach ($addr in @('D2', 'D3', 'D4', 'D5', 'D6', 'D7', 'D8', 'D11', 'D14', 'D16', 'D17', 'D18', 'D20', 'D21', 'D22', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D34', 'D36', 'D37', 'D40', 'D41', 'D42', 'D43', 'D45', 'D47', 'D48', 'D49', 'D50', 'D51')) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '62.725.74'
$ws.Range('E2').Value = '  +1.36%  '
$ws.Range('D3').Value = '3.470.34'
$ws.Range('E3').Value = '  +1.74%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '414.12'
$ws.Range('E5').Value = '  +1.46%  '
$ws.Range('D6').Value = '129.86'
$ws.Range('E6').Value = '  +0.88%  '
$ws.Range('D7').Value = '0.629'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('E10').Value = '  +10.08%  '
$ws.Range('D11').Value = '42.56'
$ws.Range('E11').Value = '  -0.22%  '
$ws.Range('E12').Value = '  +6.25%  '
$ws.Range('E13').Value = '  +3.55%  '
$ws.Range('D14').Value = '4.018.53'
$ws.Range('E14').Value = '  +1.54%  '
$ws.Range('E15').Value = '  -0.14%  '
$ws.Range('D16').Value = '20.59'
$ws.Range('E16').Value = '  -2.95%  '
$ws.Range('D17').Value = '3.461.71'
$ws.Range('E17').Value = '  +1.62%  '
$ws.Range('D18').Value = '12.64'
$ws.Range('E18').Value = '  +1.42%  '
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('D20').Value = '62.686.16'
$ws.Range('E20').Value = '  +1.33%  '
$ws.Range('D21').Value = '462.73'
$ws.Range('E21').Value = '  +2.09%  '
$ws.Range('D22').Value = '90.71'
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('E23').Value = '  +2.37%  '
$ws.Range('D24').Value = '13.29'
$ws.Range('E24').Value = '  +1.71%  '
$ws.Range('D25').Value = '10.73'
$ws.Range('E25').Value = '  +15.47%  '
$ws.Range('D26').Value = '3.33'
$ws.Range('E26').Value = '  +1.20%  '
$ws.Range('D27').Value = '33.36'
$ws.Range('E27').Value = '  +1.19%  '
$ws.Range('D28').Value = '4.80'
$ws.Range('E28').Value = '  +0.43%  '
$ws.Range('D29').Value = '7.56'
$ws.Range('E29').Value = '  -0.44%  '
$ws.Range('D30').Value = '12.08'
$ws.Range('E30').Value = '  +0.52%  '
$ws.Range('E31').Value = '  -0.96%  '
$ws.Range('E32').Value = '  -1.99%  '
$ws.Range('E33').Value = '  -1.01%  '
$ws.Range('D34').Value = '40.86'
$ws.Range('E34').Value = '  -4.34%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').Value = '58.37'
$ws.Range('E36').Value = '  +8.31%  '
$ws.Range('D37').Value = '0.0493'
$ws.Range('E37').Value = '  -1.56%  '
$ws.Range('E38').Value = '  +5.07%  '
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('D40').Value = '149.20'
$ws.Range('E40').Value = '  +4.76%  '
$ws.Range('D41').Value = '0.323'
$ws.Range('E41').Value = '  +1.80%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').Value = '0.134'
$ws.Range('E42').Value = '  +0.22%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = '2.71'
$ws.Range('E43').Value = '  +6.49%  '
$ws.Range('D45').Value = '4.40'
$ws.Range('E45').Value = '  +3.54%  '
$ws.Range('E46').Value = '  +3.80%  '
$ws.Range('B47').Value = 'PEPE'
$ws.Range('C47').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D47').Value = '0.0₃0570'
$ws.Range('E47').Value = '  +35.98%  '
$ws.Range('B48').Value = 'ThetaToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D48').Value = '2.39'
$ws.Range('E48').Value = '  +12.14%  '
$ws.Range('D49').Value = '16.44'
$ws.Range('E49').Value = '  -0.62%  '
$ws.Range('D50').Value = '22.25'
$ws.Range('E50').Value = '  -0.19%  '
$ws.Range('D51').Value = '0.141'
$ws.Range('E51').Value = '  -1.43%  '
